$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above row 159. This shifts the existing rows
# 159..193 down to 160..194, preserving all of their data untouched.
$ws.Rows.Item(159).Insert()

# Populate the new row 159 with the new record's data.
$ws.Range("A159").Value = 11
$ws.Range("B159").Value = "Vega Monumental Concepción"
$ws.Range("C159").Value = "Bíobío"
$ws.Range("D159").Value = 44782
$ws.Range("E159").Value = 8
$ws.Range("F159").Value = "Fruta"
$ws.Range("G159").Value = 100108
$ws.Range("H159").Value = "Tropicales y subtropicales"
$ws.Range("I159").Value = 100108005
$ws.Range("J159").Value = "Piña"
$ws.Range("K159").Value = "Sin especificar"
$ws.Range("L159").Value = "Segunda"
$ws.Range("M159").Value = 270
$ws.Range("N159").Value = 18000
$ws.Range("O159").Value = 19000
$ws.Range("P159").Value = 18444
$ws.Range("Q159").Value = "$/caja 14 unidades"
$ws.Range("R159").Value = "Ecuador"
$ws.Range("S159").Value = 1317
$ws.Range("T159").Value = 14
